# Update the cryptos price list (scraper refresh).
# Only cells whose text actually changed are touched; values in column D
# that would otherwise be auto-parsed by Excel as a number (e.g. "1.001",
# "286.27") are written with a leading apostrophe so they stay plain text,
# matching the workbook's original inline-string representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.408.21"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.566.46"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'286.27"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3284"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'46.59"
$ws.Range("E9").Value = "  -4.16%  "
$ws.Range("D10").Value = "'1.147"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "'0.07422"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'20.48"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'5.845"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "'6.813"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "1.588.73"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "'0.06707"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'86.24"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'6.332"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "'11.81"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "22.397.53"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'2.322"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("D26").Value = "'2.573"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "'150.76"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "'19.39"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'4.953"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "'123.70"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "1.758.29"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'1.054"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'1.975"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "'5.987"
$ws.Range("E34").Value = "  -1.97%  "
$ws.Range("D35").Value = "'9.690"
$ws.Range("E35").Value = "  -1.34%  "
$ws.Range("D36").Value = "'0.08264"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'0.02402"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "'1.306"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").Value = "'0.06325"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").Value = "'0.2187"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").Value = "'5.221"
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("D42").Value = "'11.16"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "'0.6123"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'13.67"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5956"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.748"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'2.015"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "'123.96"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'1.184"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'0.07160"
$ws.Range("E51").Value = "  -0.67%  "
